# Apply "fixed catch22 code quality" edit:
#  - Rename the catch22 feature headers (row 1) to their short/friendly names
#  - Column V ("FC_LocalSimple_mean3_stderr" -> now duplicated "periodicity")
#    gets re-populated with the values from column R (the real periodicity
#    feature) and highlighted with a solid red fill to flag the duplication.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename headers in row 1 -------------------------------------------
$headers = @(
    "mode_5",
    "mode_10",
    "stretch_high",
    "outlier_timing_pos",
    "outlier_timing_neg",
    "acf_timescale",
    "acf_first_min",
    "centroid_freq",
    "low_freq_power",
    "forecast_error",
    "trev",
    "ami2",
    "ami_timescale",
    "high_fluctuation",
    "stretch_decreasing",
    "entropy_pairs",
    "whiten_timescale",
    "periodicity",
    "dfa",
    "rs_range",
    "transition_matrix",
    "periodicity"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $headers[$i]
}

# --- 2. Recompute column V data values (copy from column R) ---------------
$ws.Range("V2:V68").Value2 = $ws.Range("R2:R68").Value2

# --- 3. Highlight the (now duplicated) column V data cells with red fill --
$ws.Range("V2:V68").Interior.Color = 255

Write-Host "catch22 header cleanup + column V fix applied"
